$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Taxonsorteringsordning (column B) values for specific rows
$ws.Range("B2").Value = 79243
$ws.Range("B3").Value = 79243
$ws.Range("B4").Value = 79243
$ws.Range("B5").Value = 79243
$ws.Range("B6").Value = 79243
$ws.Range("B8").Value = 57884
$ws.Range("B9").Value = 79243
$ws.Range("B10").Value = 79243
$ws.Range("B11").Value = 79243
$ws.Range("B12").Value = 79243
$ws.Range("B13").Value = 79243
$ws.Range("B14").Value = 57884
$ws.Range("B15").Value = 79243
